$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the summary header "Total" (D1) to "Total Per Room"
$ws.Range("D1").Value = "Total Per Room"

# Clear the Avg/Night value for "Leonardo Royal London St Paul's" (row 3),
# restoring the cell's format back to the plain (non-currency) look used by
# the rest of that column so only the borders remain.
$ws.Range("B3").ClearContents()
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# Widen column A slightly to fit the new content
$ws.Columns("A").ColumnWidth = 31.5

# Move the active selection, matching where the user left off working
$ws.Range("G9").Select()
